$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = 'Última actualización: 11:54:22'
$ws.Range("A3").Value = 'Total filas: 215'
$ws.Range("C32").Value = '215C_EL PATO'
$ws.Range("C33").Value = '14_ABASTO'
$ws.Range("C52").Value = '11_ETCHEVERRY'
$ws.Range("A53").Value = '05:55:02'
$ws.Range("C53").Value = '16_SANTA ANA'
$ws.Range("D53").Value = 97
$ws.Range("A54").Value = '06:54:06'
$ws.Range("C54").Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Range("D54").Value = 38
$ws.Range("A71").Value = '08:01:10'
$ws.Range("C71").Value = '16_SANTA ANA'
$ws.Range("D71").Value = 10
$ws.Range("A72").Value = '07:17:59'
$ws.Range("C72").Value = '15_ABASTO'
$ws.Range("D72").Value = 54
$ws.Range("A97").Value = '07:17:59'
$ws.Range("C97").Value = '17_ROMERO'
$ws.Range("D97").Value = 96
$ws.Range("A98").Value = '08:47:26'
$ws.Range("C98").Value = '10_OLMOS'
$ws.Range("D98").Value = 6
$ws.Range("A114").Value = '08:55:01'
$ws.Range("C114").Value = '16_SANTA ANA'
$ws.Range("D114").Value = 28
$ws.Range("A115").Value = '08:31:01'
$ws.Range("C115").Value = '17_ROMERO'
$ws.Range("D115").Value = 52
$ws.Range("A116").Value = '08:47:26'
$ws.Range("C116").Value = '11_ETCHEVERRY'
$ws.Range("D116").Value = 36
$ws.Range("A143").Value = '08:55:01'
$ws.Range("C143").Value = '10_OLMOS'
$ws.Range("D143").Value = 91
$ws.Range("A144").Value = '10:25:24'
$ws.Range("C144").Value = '215A_EL PATO'
$ws.Range("D144").Value = 1
$ws.Range("C170").Value = '86_EST CHICA-ESC AGRARIA'
$ws.Range("C171").Value = '16_SANTA ANA'
$ws.Range("C174").Value = '16_SANTA ANA'
$ws.Range("C175").Value = '15_ABASTO'
$ws.Range("A187").Value = '11:54:22'
$ws.Range("D187").Value = 5
$ws.Range("A189").Value = '11:54:22'
$ws.Range("D189").Value = 8
$ws.Range("A190").Value = '11:54:22'
$ws.Range("B190").Value = '12:05'
$ws.Range("C190").Value = '23_HERNANDEZ'
$ws.Range("D190").Value = 11
$ws.Range("A191").Value = '11:54:22'
$ws.Range("C191").Value = '16_P MOR-SANTA ANA'
$ws.Range("D191").Value = 12
$ws.Range("A192").Value = '11:54:22'
$ws.Range("B192").Value = '12:06'
$ws.Range("C192").Value = '14_ABASTO'
$ws.Range("D192").Value = 12
$ws.Range("C193").Value = '16_P MOR-SANTA ANA'
$ws.Range("A194").Value = '10:58:14'
$ws.Range("B194").Value = '12:07'
$ws.Range("C194").Value = '14_ABASTO'
$ws.Range("D194").Value = 69
$ws.Range("A195").Value = '11:54:22'
$ws.Range("B195").Value = '12:14'
$ws.Range("C195").Value = '17_ROMERO'
$ws.Range("D195").Value = 20
$ws.Range("A196").Value = '11:54:22'
$ws.Range("B196").Value = '12:17'
$ws.Range("C196").Value = '16_SANTA ANA'
$ws.Range("D196").Value = 23
$ws.Range("A197").Value = '11:27:22'
$ws.Range("B197").Value = '12:20'
$ws.Range("C197").Value = '14_ABASTO'
$ws.Range("D197").Value = 53
$ws.Range("A198").Value = '11:54:22'
$ws.Range("B198").Value = '12:20'
$ws.Range("C198").Value = '215A_EL PATO'
$ws.Range("D198").Value = 26
$ws.Range("C199").Value = '215A_EL PATO'
$ws.Range("A200").Value = '10:58:14'
$ws.Range("B200").Value = '12:21'
$ws.Range("C200").Value = '14_ABASTO'
$ws.Range("D200").Value = 83
$ws.Range("A201").Value = '11:54:22'
$ws.Range("B201").Value = '12:21'
$ws.Range("C201").Value = '26_HERNANDEZ'
$ws.Range("D201").Value = 27
$ws.Range("B202").Value = '12:34'
$ws.Range("C202").Value = '23_HERNANDEZ'
$ws.Range("D202").Value = 67
$ws.Range("A203").Value = '11:54:22'
$ws.Range("B203").Value = '12:35'
$ws.Range("C203").Value = '23_HERNANDEZ'
$ws.Range("D203").Value = 41
$ws.Range("A204").Value = '11:54:22'
$ws.Range("B204").Value = '12:37'
$ws.Range("C204").Value = '27_EL RETIRO'
$ws.Range("D204").Value = 43
$ws.Range("A205").Value = '11:54:22'
$ws.Range("B205").Value = '12:38'
$ws.Range("C205").Value = '17_179 Y 38'
$ws.Range("D205").Value = 44
$ws.Range("A206").Value = '11:54:22'
$ws.Range("B206").Value = '12:38'
$ws.Range("C206").Value = '11_ETCHEVERRY'
$ws.Range("D206").Value = 44
$ws.Range("A207").Value = '11:54:22'
$ws.Range("B207").Value = '12:41'
$ws.Range("C207").Value = '10_OLMOS'
$ws.Range("D207").Value = 47
$ws.Range("A208").Value = '11:54:22'
$ws.Range("B208").Value = '12:44'
$ws.Range("C208").Value = '16_SANTA ANA'
$ws.Range("D208").Value = 50
$ws.Range("A209").Value = '11:54:22'
$ws.Range("B209").Value = '12:48'
$ws.Range("C209").Value = '11_ETCHEVERRY'
$ws.Range("D209").Value = 54
$ws.Range("A210").Value = '11:54:22'
$ws.Range("B210").Value = '12:55'
$ws.Range("C210").Value = '10_OLMOS'
$ws.Range("D210").Value = 61
$ws.Range("E210").Value = 'LP1912'
$ws.Range("A211").Value = '11:54:22'
$ws.Range("B211").Value = '13:02'
$ws.Range("C211").Value = '15_ABASTO'
$ws.Range("D211").Value = 68
$ws.Range("E211").Value = 'LP1912'
$ws.Range("A212").Value = '11:54:22'
$ws.Range("B212").Value = '13:04'
$ws.Range("C212").Value = '10_OLMOS'
$ws.Range("D212").Value = 70
$ws.Range("E212").Value = 'LP1912'
$ws.Range("A213").Value = '11:54:22'
$ws.Range("B213").Value = '13:06'
$ws.Range("C213").Value = '16_P MOR-SANTA ANA'
$ws.Range("D213").Value = 72
$ws.Range("E213").Value = 'LP1912'
$ws.Range("A214").Value = '11:54:22'
$ws.Range("B214").Value = '13:14'
$ws.Range("C214").Value = '215D_EL PATO'
$ws.Range("D214").Value = 80
$ws.Range("E214").Value = 'LP1912'
$ws.Range("A215").Value = '11:54:22'
$ws.Range("B215").Value = '13:19'
$ws.Range("C215").Value = '10_OLMOS'
$ws.Range("D215").Value = 85
$ws.Range("E215").Value = 'LP1912'
$ws.Range("A216").Value = '11:54:22'
$ws.Range("B216").Value = '13:21'
$ws.Range("C216").Value = '26_HERNANDEZ'
$ws.Range("D216").Value = 87
$ws.Range("E216").Value = 'LP1912'
$ws.Range("A217").Value = '11:54:22'
$ws.Range("B217").Value = '13:26'
$ws.Range("C217").Value = '15_ABASTO'
$ws.Range("D217").Value = 92
$ws.Range("E217").Value = 'LP1912'
$ws.Range("A218").Value = '11:54:22'
$ws.Range("B218").Value = '13:26'
$ws.Range("C218").Value = '14_ABASTO'
$ws.Range("D218").Value = 92
$ws.Range("E218").Value = 'LP1912'
$ws.Range("A219").Value = '11:54:22'
$ws.Range("B219").Value = '13:46'
$ws.Range("C219").Value = '17_ROMERO'
$ws.Range("D219").Value = 112
$ws.Range("E219").Value = 'LP1912'
$ws.Range("A220").Value = '11:54:22'
$ws.Range("B220").Value = '13:50'
$ws.Range("C220").Value = '215A_EL PATO'
$ws.Range("D220").Value = 116
$ws.Range("E220").Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = 'Última actualización: 11:54:22'
$ws.Range("A3").Value = 'Total filas: 26'
$ws.Range("A28").Value = '11:54:22'
$ws.Range("D28").Value = 26
$ws.Range("A30").Value = '11:54:22'
$ws.Range("D30").Value = 80
$ws.Range("E30").Value = 'LP1912'
$ws.Range("A31").Value = '11:54:22'
$ws.Range("B31").Value = '13:50'
$ws.Range("C31").Value = '215A_EL PATO'
$ws.Range("D31").Value = 116
$ws.Range("E31").Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = 'Última actualización: 11:54:22'
$ws.Range("A3").Value = 'Total filas: 26'
$ws.Range("A29").Value = '11:54:22'
$ws.Range("D29").Value = 10
$ws.Range("A30").Value = '11:54:22'
$ws.Range("D30").Value = 60
$ws.Range("E30").Value = 'L6203'
$ws.Range("A31").Value = '11:54:22'
$ws.Range("B31").Value = '13:31'
$ws.Range("C31").Value = '215B_LP-P MOR-1 Y 57'
$ws.Range("D31").Value = 97
$ws.Range("E31").Value = 'L6173'
